# Update betting odds figures on the weekly FlashScore games sheet.
# Each entry maps a cell address to its new numeric value, matching the
# changes described in the source diff (rows 6, 7, 8, 9, 13 and 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    # Row 6
    "G6"  = 2.05
    "H6"  = 3.3
    "I6"  = 3.7
    "J6"  = 1.08
    "K6"  = 8
    "U6"  = 9
    "V6"  = 9
    "W6"  = 17
    "AE6" = 9.5
    "AI6" = 34

    # Row 7
    "G7" = 1.73
    "P7" = 1.44
    "Q7" = 2.63

    # Row 8
    "L8"  = 1.31
    "M8"  = 2.87
    "N8"  = 1.91
    "O8"  = 1.7
    "Q8"  = 2.42
    "R8"  = 1.7
    "S8"  = 1.93
    "T8"  = 9.25
    "U8"  = 16
    "V8"  = 10.25
    "X8"  = 25
    "Y8"  = 30
    "Z8"  = 8.75
    "AB8" = 13
    "AC8" = 60
    "AD8" = 450
    "AF8" = 11.5
    "AH8" = 26

    # Row 9
    "G9"  = 3.7
    "H9"  = 3.35
    "I9"  = 1.91
    "R9"  = 1.98
    "S9"  = 1.65
    "T9"  = 9
    "U9"  = 18.5
    "V9"  = 13.5
    "W9"  = 55
    "X9"  = 40
    "Y9"  = 55
    "AA9" = 6.5
    "AB9" = 18
    "AC9" = 110
    "AE9" = 6
    "AF9" = 8.25
    "AH9" = 16
    "AI9" = 17.5
    "AJ9" = 35

    # Row 13
    "T13"  = 8
    "W13"  = 23
    "AI13" = 23
    "AJ13" = 34

    # Row 14
    "G14"  = 1.57
    "H14"  = 3.85
    "I14"  = 4.85
    "J14"  = 1.05
    "K14"  = 7.7
    "L14"  = 1.27
    "M14"  = 3.4
    "N14"  = 1.82
    "O14"  = 1.9
    "P14"  = 1.37
    "Q14"  = 2.85
    "R14"  = 1.87
    "S14"  = 1.83
    "T14"  = 6.8
    "U14"  = 7.3
    "W14"  = 11.5
    "X14"  = 13
    "Y14"  = 27
    "Z14"  = 7.7
    "AA14" = 7.6
    "AB14" = 17.5
    "AE14" = 13.5
    "AF14" = 28
    "AG14" = 16
    "AH14" = 90
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
